$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.232.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4696"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06569"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.80"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07992"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.70"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.862.03"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.109"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6784"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.209.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.59"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007645"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.105.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.189"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.19"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.189"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.370"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09910"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.346"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04711"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7031"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.704"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01874"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.607"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.336"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.48"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.939"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8396"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.77"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9985"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.060"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.149"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "934.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.09"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.43%  "
